# Automatische test-sync: 2025-08-14 20:24:50
# Adds a new log row (row 4) to the "Logs" sheet, extends the conditional
# formatting ranges that cover the log table, and bumps the "Aantal" count
# on the "Dashboard" sheet for the matching category.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row of data
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Demo inplannen"
$logs.Range("B4").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C4").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D4").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E4").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Range("F4").Value = "2025-08-14 20:24:08"
$logs.Range("G4").Value = "Nee"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting ranges so row 4 is included too
#    (D2:D3 -> D2:D4, G2:G3 -> G2:G4, H2:H3 -> H2:H4, I2:I3 -> I2:I4,
#     J2:J3 -> J2:J4)
# ---------------------------------------------------------------------
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range($col + "2:" + $col + "3")
    $newRange = $logs.Range($col + "2:" + $col + "4")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count(); $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: bump the count for "Intern verzoek / Actie voor
#    medewerker" from 2 to 3
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 3
